# Correccion a relation field grande
#
# Sheet1 (create-fields): insert two new columns -
#   - "relation"      before the "multiple" column (new column G)
#   - "relation_cond" before the "tooltip" column (new column N, after first insert)
# Sheet2 (edit-fields): rework the field-metadata rows - drop the "barcode"
#   node rows, add "filter"/"page" rows at the top and "relation" rows for
#   filter/activity/indicator nodes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: create-fields
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("create-fields")

# Insert new column G ("relation"), shifting old G..O to H..P
$ws1.Columns.Item(7).Insert()
$ws1.Cells.Item(1, 7).Value2 = "relation"
for ($r = 2; $r -le 7; $r++) {
  $ws1.Cells.Item($r, 7).HorizontalAlignment = -4131
  $ws1.Cells.Item($r, 7).VerticalAlignment = -4160
}
$ws1.Columns.Item(7).ColumnWidth = 8

# Insert new column N ("relation_cond"), shifting old N..O (now at N..O) right
$ws1.Columns.Item(14).Insert()
$ws1.Cells.Item(1, 14).Value2 = "relation_cond"
$ws1.Columns.Item(14).ColumnWidth = 13.28

$ws1.Range("N2").Select()

# ---------------------------------------------------------------------
# Sheet2: edit-fields
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("edit-fields")

# Grow the table from 43 to 46 data rows (3 new rows) before rewriting.
$ws2.Rows.Item(44).EntireRow.Insert()
$ws2.Rows.Item(45).EntireRow.Insert()
$ws2.Rows.Item(46).EntireRow.Insert()

$ws2.Range("A2:D46").ClearContents()

$data = @(
  @("filter",          "user_id",                 "type",         "string"),
  @("filter",          "user_id",                 "relation",     0),
  @("filter",          "category_id",             "type",         "string"),
  @("filter",          "category_id",             "relation",     0),
  @("page",            "name",                     "display_list", "show"),
  @("role",            "permission_role",          "multiple",     1),
  @("user",            "indicator_alert_users",    "multiple",     1),
  @("user",            "indicator_graph_users",    "multiple",     1),
  @("user",            "email",                    "required",     0),
  @("user",            "cellphone",                "required",     0),
  @("user",            "username",                 "required",     0),
  @("user",            "status",                   "type",         "radio"),
  @("user",            "notifications_email",      "type",         "radio"),
  @("user",            "notifications_sms",        "type",         "radio"),
  @("user",            "notifications_app",        "type",         "radio"),
  @("user",            "remember_token",           "display_item", "none"),
  @("user",            "notifications_last_read",  "display_item", "none"),
  @("user",            "last_activity",            "display_item", "none"),
  @("user",            "last_login",               "display_item", "none"),
  @("user",            "last_session",             "display_item", "none"),
  @("user",            "role_user",                "required",     1),
  @("menu",            "level",                    "preset",       1),
  @("menu",            "order",                    "display_list", "none"),
  @("menu",            "parent_id",                "display_list", "none"),
  @("menu",            "menu_type",                "display_list", "show"),
  @("menu",            "order",                    "display_list", "show"),
  @("menu",            "name",                     "display_list", "show"),
  @("alert",           "node_id",                  "preset",       1),
  @("activity",        "item_id",                  "type",         "string"),
  @("activity",        "item_id",                  "relation",     0),
  @("inbox",           "from_user_id",              "value",        "user"),
  @("inbox",           "to_user_id",                "value",        "user"),
  @("variable",        "value",                     "display_list", "show"),
  @("indicator",       "node_id",                   "preset",       1),
  @("indicator",       "user_id",                   "type",         "string"),
  @("indicator",       "user_id",                   "relation",     0),
  @("indicator",       "user_id",                   "display_list", "none"),
  @("indicator",       "indicator_alerts",          "display_list", "excel"),
  @("indicator",       "indicator_graphs",          "display_list", "excel"),
  @("indicator",       "indicator_values",          "display_list", "excel"),
  @("indicator",       "user_id",                   "display_item", "none"),
  @("indicator",       "formula",                   "message",      "Escriba una formula bajo lo indicado en el manual, solo para uso avanzado."),
  @("indicator",       "result_custom",             "message",      "Utilice solo si el manual lo indica y siga las intrucciones."),
  @("indicator-graph", "name",                      "display_item", "none"),
  @("indicator-alert", "name",                      "display_item", "none")
)

$r = 2
foreach ($row in $data) {
  $ws2.Cells.Item($r, 1).Value2 = $row[0]
  $ws2.Cells.Item($r, 2).Value2 = $row[1]
  $ws2.Cells.Item($r, 3).Value2 = $row[2]
  $ws2.Cells.Item($r, 4).Value2 = $row[3]
  $r = $r + 1
}

$ws2.Range("C16").Select()
